$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.346.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.090.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  -6.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.366'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.089.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.744'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.912.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.661.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.128.59'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("B19").Value = 'SuiNetwork'
$ws.Range("C19").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000223'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '435.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.17%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.199'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.88%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.81'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.74%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.152'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '503.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0877'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.405'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +54.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.691'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '152.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.73%  '
